# Auto update Excel log
# Appends newly-logged sensor events (2026-01-28, ~12:35-12:36) to the
# PIR, Humidity, Temperature, Proximity and Camera sheets.

$wb = $excel.ActiveWorkbook

# Helper: write literal-text values into a range without Excel's COM
# layer auto-coercing date-looking ("YYYY-MM-DD") or percent-looking
# ("NN.N%") strings into numeric serials. Applying a Text number format
# before the assignment forces a string; ClearFormats() afterwards
# removes the now-unneeded style so the cell matches the workbook's
# plain default styling (same as every other text cell in these logs).
function Set-TextValue($Range, $Value) {
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

# ---------------------------------------------------------------------
# PIR sheet: rows 373-385 (Bathroom / No Motion / Inactive)
# ---------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirStart = 373
$pirTimes = @("12:35:32","12:35:37","12:35:43","12:35:46","12:35:49","12:35:56","12:36:00","12:36:04","12:36:09","12:36:14","12:36:19","12:36:26","12:36:30")
$pirCount = $pirTimes.Count
$pirEnd = $pirStart + $pirCount - 1

Set-TextValue $wsPIR.Range("A$($pirStart):A$($pirEnd)") "2026-01-28"

$pirData = New-Object 'object[,]' $pirCount,5
for ($i = 0; $i -lt $pirCount; $i++) {
    $pirData[$i,0] = $pirTimes[$i]
    $pirData[$i,1] = "12:00"
    $pirData[$i,2] = "Bathroom"
    $pirData[$i,3] = "No Motion"
    $pirData[$i,4] = "Inactive"
}
$wsPIR.Range("B$($pirStart):F$($pirEnd)").Value = $pirData

# ---------------------------------------------------------------------
# Humidity sheet: rows 349-361 (Bathroom / xx.x% / Active)
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humStart = 349
$humTimes = @("12:35:33","12:35:35","12:35:41","12:35:44","12:35:47","12:35:50","12:35:54","12:35:58","12:36:02","12:36:06","12:36:10","12:36:24","12:36:31")
$humValues = @("87.1%","86.1%","87.0%","86.2%","87.1%","87.1%","86.1%","87.1%","86.1%","87.1%","87.0%","86.1%","87.0%")
$humCount = $humTimes.Count
$humEnd = $humStart + $humCount - 1

Set-TextValue $wsHumidity.Range("A$($humStart):A$($humEnd)") "2026-01-28"

$humDataBCD = New-Object 'object[,]' $humCount,3
for ($i = 0; $i -lt $humCount; $i++) {
    $humDataBCD[$i,0] = $humTimes[$i]
    $humDataBCD[$i,1] = "12:00"
    $humDataBCD[$i,2] = "Bathroom"
}
$wsHumidity.Range("B$($humStart):D$($humEnd)").Value = $humDataBCD

# Column E (percent readings) must stay literal text, not numeric %.
for ($i = 0; $i -lt $humCount; $i++) {
    $rowNum = $humStart + $i
    Set-TextValue $wsHumidity.Range("E$($rowNum)") $humValues[$i]
}

$humDataF = New-Object 'object[,]' $humCount,1
for ($i = 0; $i -lt $humCount; $i++) { $humDataF[$i,0] = "Active" }
$wsHumidity.Range("F$($humStart):F$($humEnd)").Value = $humDataF

# ---------------------------------------------------------------------
# Temperature sheet: rows 349-360 (Bathroom / 23.0C / Active)
# ---------------------------------------------------------------------
$wsTemp = $wb.Worksheets.Item("Temperature")
$tempStart = 349
$tempTimes = @("12:35:34","12:35:36","12:35:42","12:35:45","12:35:48","12:35:51","12:35:55","12:35:59","12:36:03","12:36:07","12:36:11","12:36:25")
$tempCount = $tempTimes.Count
$tempEnd = $tempStart + $tempCount - 1

Set-TextValue $wsTemp.Range("A$($tempStart):A$($tempEnd)") "2026-01-28"

$tempData = New-Object 'object[,]' $tempCount,5
for ($i = 0; $i -lt $tempCount; $i++) {
    $tempData[$i,0] = $tempTimes[$i]
    $tempData[$i,1] = "12:00"
    $tempData[$i,2] = "Bathroom"
    $tempData[$i,3] = "23.0C"
    $tempData[$i,4] = "Active"
}
$wsTemp.Range("B$($tempStart):F$($tempEnd)").Value = $tempData

# ---------------------------------------------------------------------
# Proximity sheet: rows 13-16 (Living Room Main Entrance)
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
$proxStart = 13
$proxRows = @(
    @("12:35:38","Detected","Active"),
    @("12:35:40","Clear","Inactive"),
    @("12:36:20","Detected","Active"),
    @("12:36:23","Clear","Inactive")
)
$proxCount = $proxRows.Count
$proxEnd = $proxStart + $proxCount - 1

Set-TextValue $wsProximity.Range("A$($proxStart):A$($proxEnd)") "2026-01-28"

$proxData = New-Object 'object[,]' $proxCount,5
for ($i = 0; $i -lt $proxCount; $i++) {
    $proxData[$i,0] = $proxRows[$i][0]
    $proxData[$i,1] = "12:00"
    $proxData[$i,2] = "Living Room Main Entrance"
    $proxData[$i,3] = $proxRows[$i][1]
    $proxData[$i,4] = $proxRows[$i][2]
}
$wsProximity.Range("B$($proxStart):F$($proxEnd)").Value = $proxData

# ---------------------------------------------------------------------
# Camera sheet: rows 5-6 (Living Room Main Entrance, columns A-D only)
# ---------------------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")
$camStart = 5
$camTimes = @("12:35:39","12:36:22")
$camCount = $camTimes.Count
$camEnd = $camStart + $camCount - 1

Set-TextValue $wsCamera.Range("A$($camStart):A$($camEnd)") "2026-01-28"

$camData = New-Object 'object[,]' $camCount,3
for ($i = 0; $i -lt $camCount; $i++) {
    $camData[$i,0] = $camTimes[$i]
    $camData[$i,1] = "12:00"
    $camData[$i,2] = "Living Room Main Entrance"
}
$wsCamera.Range("B$($camStart):D$($camEnd)").Value = $camData
